# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 114 (pushing the existing rows
# 114-144 down to 115-145) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 114..144 down by one row.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(114, 1).Value  = 1
$ws.Cells.Item(114, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(114, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(114, 4).Value  = 44876
$ws.Cells.Item(114, 5).Value  = 15
$ws.Cells.Item(114, 6).Value  = 100112036
$ws.Cells.Item(114, 7).Value  = "Caigua"
$ws.Cells.Item(114, 8).Value  = "Sin especificar"
$ws.Cells.Item(114, 9).Value  = "Primera"
$ws.Cells.Item(114, 10).Value = 100
$ws.Cells.Item(114, 11).Value = 9000
$ws.Cells.Item(114, 12).Value = 10000
$ws.Cells.Item(114, 13).Value = 9500
$ws.Cells.Item(114, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(114, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(114, 16).Value = 475
$ws.Cells.Item(114, 17).Value = 20
$ws.Cells.Item(114, 18).Value = "Hortaliza"

# Match the date formatting style used by the rest of column D.
$ws.Cells.Item(114, 4).NumberFormat = $ws.Cells.Item(115, 4).NumberFormat()
